# Update "want to go" (F column) counts across the four worksheets to
# reflect newly generated numbers (gh-pages output regeneration).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 262
$ws1.Range("F4").Value  = 857
$ws1.Range("F6").Value  = 427
$ws1.Range("F7").Value  = 631
$ws1.Range("F8").Value  = 230
$ws1.Range("F10").Value = 365
$ws1.Range("F11").Value = 165
$ws1.Range("F12").Value = 734
$ws1.Range("F14").Value = 1859
$ws1.Range("F15").Value = 389
$ws1.Range("F16").Value = 4387
$ws1.Range("F17").Value = 388
$ws1.Range("F18").Value = 500
$ws1.Range("F19").Value = 20
$ws1.Range("F20").Value = 65
$ws1.Range("F21").Value = 154

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value  = 24
$ws2.Range("F6").Value  = 120
$ws2.Range("F7").Value  = 488
$ws2.Range("F13").Value = 100
$ws2.Range("F14").Value = 42

# --- Sheet "本地生活" (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 5394
$ws3.Range("F3").Value = 341
$ws3.Range("F4").Value = 307

# --- Sheet "全部类型" (All Types, aggregate of the above) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 5394
$ws4.Range("F4").Value  = 341
$ws4.Range("F6").Value  = 307
$ws4.Range("F7").Value  = 262
$ws4.Range("F10").Value = 24
$ws4.Range("F11").Value = 120
$ws4.Range("F12").Value = 488
$ws4.Range("F13").Value = 857
$ws4.Range("F17").Value = 427
$ws4.Range("F18").Value = 631
$ws4.Range("F19").Value = 230
$ws4.Range("F22").Value = 365
$ws4.Range("F23").Value = 165
$ws4.Range("F26").Value = 734
$ws4.Range("F28").Value = 100
$ws4.Range("F29").Value = 1859
$ws4.Range("F30").Value = 389
$ws4.Range("F31").Value = 4387
$ws4.Range("F32").Value = 42
$ws4.Range("F33").Value = 388
$ws4.Range("F34").Value = 500
$ws4.Range("F35").Value = 20
$ws4.Range("F36").Value = 65
$ws4.Range("F38").Value = 154
